$d = $word.ActiveDocument

# Table 1 is a single 20-row x 5-col grid; only the 5 "data" rows
# (1, 5, 9, 13, 17) hold text, the rest are blank spacer rows.
# Cells are addressed by (row, column) so each value is updated in
# place regardless of text collisions between old/new values.

$d.Tables(1).Cell(1, 1).Range.Text = "75÷6=12, 3"
$d.Tables(1).Cell(1, 2).Range.Text = "42÷4=10, 2"
$d.Tables(1).Cell(1, 3).Range.Text = "58÷3=19, 1"
$d.Tables(1).Cell(1, 4).Range.Text = "50÷2=25, 0"
$d.Tables(1).Cell(1, 5).Range.Text = "44÷9=4, 8"

$d.Tables(1).Cell(5, 1).Range.Text = "26÷2=13, 0"
$d.Tables(1).Cell(5, 2).Range.Text = "47÷8=5, 7"
$d.Tables(1).Cell(5, 3).Range.Text = "88÷6=14, 4"
$d.Tables(1).Cell(5, 4).Range.Text = "42÷5=8, 2"
$d.Tables(1).Cell(5, 5).Range.Text = "53÷6=8, 5"

$d.Tables(1).Cell(9, 1).Range.Text = "23÷9=2, 5"
$d.Tables(1).Cell(9, 2).Range.Text = "63÷7=9, 0"
$d.Tables(1).Cell(9, 3).Range.Text = "19÷3=6, 1"
$d.Tables(1).Cell(9, 4).Range.Text = "62÷6=10, 2"
$d.Tables(1).Cell(9, 5).Range.Text = "74÷8=9, 2"

$d.Tables(1).Cell(13, 1).Range.Text = "61÷3=20, 1"
$d.Tables(1).Cell(13, 2).Range.Text = "50÷4=12, 2"
$d.Tables(1).Cell(13, 3).Range.Text = "45÷7=6, 3"
$d.Tables(1).Cell(13, 4).Range.Text = "24÷7=3, 3"
$d.Tables(1).Cell(13, 5).Range.Text = "31÷7=4, 3"

$d.Tables(1).Cell(17, 1).Range.Text = "34÷2=17, 0"
$d.Tables(1).Cell(17, 2).Range.Text = "99÷9=11, 0"
$d.Tables(1).Cell(17, 3).Range.Text = "87÷7=12, 3"
$d.Tables(1).Cell(17, 4).Range.Text = "90÷2=45, 0"
$d.Tables(1).Cell(17, 5).Range.Text = "34÷7=4, 6"

$d.Save()
